# The workbook currently has two sheets in this order: "2022-Q2", "总计".
# Re-sort the tabs so the summary sheet "总计" comes first, followed by
# the quarterly detail sheet "2022-Q2" (matches the author's "resort
# sheetname" commit - data/content of each sheet is left untouched).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Move "总计" so it sits immediately before "2022-Q2" -> new order:
# 总计, 2022-Q2
$wsTotal.Move($wsQ2)
